# Update weekly excess mortality analysis (CBS oversterfte workbook)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Revise a handful of previously-reported weekly death counts (column G) ---
# The formulas in column I (Waargenomen - Verwacht) recalculate automatically.
$ws.Range("G8").Value  = 4305
$ws.Range("G21").Value = 2528
$ws.Range("G23").Value = 2667
$ws.Range("G24").Value = 2640
$ws.Range("G26").Value = 2853
$ws.Range("G31").Value = 2891
$ws.Range("G33").Value = 3019
$ws.Range("G34").Value = 3212
$ws.Range("G35").Value = 3444
$ws.Range("G36").Value = 3674
$ws.Range("G37").Value = 3587
$ws.Range("G38").Value = 3552
$ws.Range("G39").Value = 3315
$ws.Range("G40").Value = 3373

# --- Insert week 49, pushing the "Totaal" row from 42 down to 43 ---
$ws.Rows.Item(41).Insert()

$ws.Range("F41").Value = 49
$ws.Range("G41").Value = 3448
$ws.Range("H41").Value = 3037
$ws.Range("I41").Formula = "=G41-H41"

# Update the selection to reflect the new layout
$ws.Range("J41").Select() | Out-Null
